$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Drop the oldest quarter (column D) and shift everything left ---
$ws.Columns("D").Delete()

# --- 2) Bring column M back to life with the same look as column L ---
$ws.Range("L8:L27").Copy()
$ws.Range("M8:M27").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Columns("M").ColumnWidth = $ws.Columns("L").ColumnWidth

# --- 3) New period header or column M ---
$ws.Range("M8").Value = "9 ماهه منتهی به 1401/12"

# --- 4) Publication-date row: two revised dates + the new quarter's date ---
$ws.Range("I9").Value = "1402-01-30 (2)"
$ws.Range("J9").Value = "1402-01-30 (7)"

$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-01-30"
$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)       # xlPasteFormats (restore style lost by NumberFormat change)
$ws.Application.CutCopyMode = $false

# --- 5) New quarter's figures ---
$ws.Range("M11").Value = 8211682
$ws.Range("M12").Value = -4343758
$ws.Range("M13").Value = 3867924
$ws.Range("M14").Value = -331097
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 57480
$ws.Range("M17").Value = 3594307
$ws.Range("M18").Value = -779182
$ws.Range("M19").Value = 1209757
$ws.Range("M20").Value = 4024882
$ws.Range("M21").Value = -393269
$ws.Range("M22").Value = 3631613
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 3631613
$ws.Range("M25").Value = 484
$ws.Range("M26").Value = 7500000
$ws.Range("M27").Value = 484
